# Applies the changes described by the commit:
# "se agregó el listar y se arregló el formulario del consultorio"
#
# Sheet "Data Medic" (1st sheet):
#   - B7: add "( al 100%)" to Sede / Área, remove trailing space after "Consultorio"
#   - F7: 0.28 -> 0.56 (% avance)
#   - Active cell selection moves from B7 to B9
#
# Sheet "Java Web Developer - Tesis" (2nd sheet):
#   - B6: add 3 new bullet lines (Herencia, Polimorfismo, Interfaz) under POO theory study
#   - Row 6 height grows from 90 to 135 (to fit the extra lines)
#   - F6: 0.25 -> 0.4 (% avance)
#   - Active cell selection moves from F20:F22 to F6

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data Medic")
$wsJava = $wb.Worksheets.Item("Java Web Developer - Tesis")

# ---- Sheet "Data Medic" ----

$newB7 = "Agregar el mantenimiento completo para la: `n" +
         "                                     a. Empresa ( al 100%)`n" +
         "                                     b. Sede ( al 100%)`n" +
         "                                     c. Consultorio`n" +
         "                                     d. Área ( al 100%)`n" +
         "                                     e. cita`n" +
         "                                     f. hora_atencion_doctor`n" +
         "                                     g. otra_especialización`n" +
         "Dar acceso a los usuarios correspondiente: Super Usuario: Empresa, Gerente: Sede y Consultorio (100%)"

$wsData.Range("B7").Value = $newB7
$wsData.Range("F7").Value = 0.56000000000000005

# ---- Sheet "Java Web Developer - Tesis" ----

$newB6 = "Java Web Developer: `n" +
         "                                   a. Sesión 1 (20%) minuto 55 del video.`n" +
         "                                   b. Estudiar POO teoría`n" +
         "                                       - Encapsulamiento 100%`n" +
         "                                       - Herencia 100%`n" +
         "                                       - Polimorfismo 100%`n" +
         "                                       - Interfaz 100%`n" +
         "                                   c. Prácticar POO con ejercicios que yo mismo proponga`n" +
         "                                       - Encapsulamiento 100%"

$wsJava.Range("B6").Value = $newB6
$wsJava.Rows.Item(6).RowHeight = 135
$wsJava.Range("F6").Value = 0.4

# ---- Update the active selections shown in each sheet view ----

$wsData.Activate()
$wsData.Range("B9").Select()

$wsJava.Activate()
$wsJava.Range("F6").Select()

# Leave "Data Medic" as the active/visible tab, matching the original workbook.
$wsData.Activate()
